# Apply the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values (column D) are plain decimal numbers
# (e.g. "310.84"). Excel's COM layer auto-coerces such literal-looking
# strings assigned through .Value into real numbers, which would change the
# cell from a text/inline-string cell (as it is in the source workbook) into
# a numeric cell. Temporarily force a Text format on just those cells so the
# assignment keeps them as text, then clear the temporary formatting again so
# the cells end up unstyled, same as in the original file.
$riskyPriceCells = @("D5", "D7", "D9", "D10", "D11", "D15", "D19", "D22", "D23", "D25", "D28", "D29", "D32", "D35", "D37", "D38", "D44", "D46", "D48", "D49", "D50")
foreach ($addr in $riskyPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.422.46'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '2.309.50'
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '310.84'
$ws.Range("E5").Value = '  +1.54%  '
$ws.Range("E6").Value = '  +4.64%  '
$ws.Range("D7").Value = '0.536'
$ws.Range("E7").Value = '  +1.46%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  +7.56%  '
$ws.Range("D10").Value = '35.71'
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").Value = '0.0814'
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("E12").Value = '  -0.59%  '
$ws.Range("E13").Value = '  +0.72%  '
$ws.Range("D14").Value = '2.667.14'
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = '14.96'
$ws.Range("E15").Value = '  +1.26%  '
$ws.Range("D16").Value = '2.307.72'
$ws.Range("E16").Value = '  +1.94%  '
$ws.Range("E17").Value = '  +1.95%  '
$ws.Range("D18").Value = '43.349.52'
$ws.Range("E18").Value = '  +2.87%  '
$ws.Range("D19").Value = '12.32'
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").Value = '0.0₃0929'
$ws.Range("E20").Value = '  +2.52%  '
$ws.Range("E21").Value = '  +2.31%  '
$ws.Range("D22").Value = '68.05'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = '241.34'
$ws.Range("E23").Value = '  +1.36%  '
$ws.Range("E24").Value = '  +2.73%  '
$ws.Range("D25").Value = '2.62'
$ws.Range("E25").Value = '  +1.67%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -1.61%  '
$ws.Range("D28").Value = '24.61'
$ws.Range("E28").Value = '  +4.39%  '
$ws.Range("D29").Value = '36.74'
$ws.Range("E29").Value = '  -3.49%  '
$ws.Range("E30").Value = '  +1.06%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("D32").Value = '167.60'
$ws.Range("E32").Value = '  +3.99%  '
$ws.Range("E33").Value = '  +0.74%  '
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '0.0744'
$ws.Range("E35").Value = '  +0.80%  '
$ws.Range("E36").Value = '  +5.68%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").Value = '17.61'
$ws.Range("E37").Value = '  -0.56%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").Value = '3.06'
$ws.Range("E38").Value = '  -2.85%  '
$ws.Range("E39").Value = '  +3.07%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("E41").Value = '  +1.35%  '
$ws.Range("E42").Value = '  +7.01%  '
$ws.Range("E43").Value = '  -0.91%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0289'
$ws.Range("E44").Value = '  +2.64%  '
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.967.90'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '19.20'
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("E47").Value = '  +2.07%  '
$ws.Range("D48").Value = '9.93'
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("D49").Value = '55.53'
$ws.Range("E49").Value = '  +3.50%  '
$ws.Range("D50").Value = '2.92'
$ws.Range("E50").Value = '  +4.62%  '
$ws.Range("E51").Value = '  +6.56%  '

foreach ($addr in $riskyPriceCells) {
    $ws.Range($addr).ClearFormats()
}

